$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at C and D, shifting old C (full), D (tipo), E (link) to E, F, G
$ws.Columns("C:D").Insert()

# Header row
$ws.Range("C1").Value() = "modelo"
$ws.Range("D1").Value() = "politica"

# Per-row data: modelo (C), politica (D), tipo (F, lowercase), link (G)
$modelo = @{}
$politica = @{}
$tipo = @{}
$link = @{}

$modelo[2] = "Modelo identificado mas fora do range de preco"
$tipo[2] = "premium"
$link[2] = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:251253029#searchVariation=MLB21392652&position=2&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[3] = "Modelo identificado mas fora do range de preco"
$tipo[3] = "premium"
$link[3] = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:251253029#searchVariation=MLB24154371&position=3&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[4] = "FONTE 120 BOB"
$politica[4] = "Acima"
$tipo[4] = "classico"
$link[4] = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:251253029#searchVariation=MLB22144397&position=5&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[5] = "Modelo identificado mas fora do range de preco"
$tipo[5] = "classico"
$link[5] = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-cor-preto/p/MLB23998473?pdp_filters=seller_id:251253029#searchVariation=MLB23998473&position=6&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[6] = "FONTE 200 BOB"
$politica[6] = "Igual"
$tipo[6] = "classico"
$link[6] = "https://produto.mercadolivre.com.br/MLB-4050176664-fonte-automotiva-bivolt-jfa-bob-storm-200a-dinamico-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[7] = "Sem Modelo"
$tipo[7] = "classico"
$link[7] = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:251253029#searchVariation=MLB27687422&position=7&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[8] = "Modelo identificado mas fora do range de preco"
$tipo[8] = "premium"
$link[8] = "https://produto.mercadolivre.com.br/MLB-3487253887-fonte-carregador-de-bateria-jfa-60a-lite-storm-slim-bivolt-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[9] = "Sem Modelo"
$tipo[9] = "classico"
$link[9] = "https://www.mercadolivre.com.br/conversor-fio-para-rca-remoto-slim-12v-jfa-automotivo-cd-dvd/p/MLB25707531?pdp_filters=seller_id:251253029#searchVariation=MLB25707531&position=1&search_layout=stack&type=product&tracking_id=9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[10] = "Sem Modelo"
$tipo[10] = "classico"
$link[10] = "https://produto.mercadolivre.com.br/MLB-4050015720-controle-jfa-longa-distancia-k1200-preto-com-verde-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[11] = "FONTE 120 BOB"
$politica[11] = "Acima"
$tipo[11] = "premium"
$link[11] = "https://produto.mercadolivre.com.br/MLB-3450034619-fonte-carregador-jfa-120a-bob-storm-bivolt-automatico-_JM#position%3D11%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[12] = "Modelo identificado mas fora do range de preco"
$tipo[12] = "classico"
$link[12] = "https://produto.mercadolivre.com.br/MLB-3450026499-fonte-jfa-automotivo-200a-storm-lite-12v-bivolt-azul-_JM#position%3D12%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[13] = "Modelo identificado mas fora do range de preco"
$tipo[13] = "classico"
$link[13] = "https://produto.mercadolivre.com.br/MLB-3487170605-fonte-e-carregador-de-bateria-jfa-60a-lite-storm-slim-bivolt-_JM#position%3D13%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[14] = "Sem Modelo"
$tipo[14] = "premium"
$link[14] = "https://produto.mercadolivre.com.br/MLB-3449829257-controle-longa-distancia-jfa-acqua-branco-resistente-a-agua-_JM#position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[15] = "Sem Modelo"
$tipo[15] = "premium"
$link[15] = "https://produto.mercadolivre.com.br/MLB-3449675865-controle-longa-distancia-k1200-preto-1200m-jfa-eletronico-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[16] = "Sem Modelo"
$tipo[16] = "premium"
$link[16] = "https://produto.mercadolivre.com.br/MLB-3449630269-controle-longa-distancia-k1200-azul-1200m-jfa-eletronico-_JM#position%3D16%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[17] = "Sem Modelo"
$tipo[17] = "classico"
$link[17] = "https://produto.mercadolivre.com.br/MLB-4050100920-voltimetro-sequenciador-jfa-vs5hi-3-em-1-voltagem-12v-_JM#position%3D17%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[18] = "Sem Modelo"
$tipo[18] = "premium"
$link[18] = "https://produto.mercadolivre.com.br/MLB-3449929417-filtro-anti-ruido-jfa-rca-eletromagnetico-stereo-_JM#position%3D18%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[19] = "Sem Modelo"
$tipo[19] = "premium"
$link[19] = "https://produto.mercadolivre.com.br/MLB-4049938564-controle-longa-distancia-jfa-k1200-preto-com-verde-_JM#position%3D19%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[20] = "Sem Modelo"
$tipo[20] = "premium"
$link[20] = "https://produto.mercadolivre.com.br/MLB-3449904347-controle-jfa-longa-distancia-redline-1200m-entrada-wr-_JM#position%3D20%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[21] = "Sem Modelo"
$tipo[21] = "premium"
$link[21] = "https://produto.mercadolivre.com.br/MLB-4050053906-controle-longa-distancia-jfa-acqua-preto-resistente-a-agua-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[22] = "Sem Modelo"
$tipo[22] = "classico"
$link[22] = "https://produto.mercadolivre.com.br/MLB-3449818587-controle-jfa-longa-distancia-k1200-vermelho-1200m-_JM#position%3D22%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[23] = "Sem Modelo"
$tipo[23] = "classico"
$link[23] = "https://produto.mercadolivre.com.br/MLB-4049977502-controle-jfa-longa-distancia-k1200-azul-1200m-_JM#position%3D23%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[24] = "Sem Modelo"
$tipo[24] = "classico"
$link[24] = "https://produto.mercadolivre.com.br/MLB-4050027144-controle-jfa-longa-distancia-k1200-preto-com-laranja-_JM#position%3D24%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[25] = "Sem Modelo"
$tipo[25] = "premium"
$link[25] = "https://produto.mercadolivre.com.br/MLB-3449802721-controle-longa-distancia-jfa-k1200-preto-com-laranja-_JM#position%3D25%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[26] = "Sem Modelo"
$tipo[26] = "classico"
$link[26] = "https://produto.mercadolivre.com.br/MLB-4049954464-controle-jfa-longa-distancia-k1200-preto-1200m-_JM#position%3D26%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[27] = "Sem Modelo"
$tipo[27] = "premium"
$link[27] = "https://produto.mercadolivre.com.br/MLB-3449889659-voltimetro-sequenciador-jfa-vs5hi-3-em-1-voltagem-12v-_JM#position%3D27%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

$modelo[28] = "FONTE 120A"
$politica[28] = "Igual"
$tipo[28] = "classico"
$link[28] = "https://produto.mercadolivre.com.br/MLB-2731042154-fonte-automotiva-120a-storm-jfa-carregador-_JM#position%3D28%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D9e2739fe-416b-492a-88fb-a96f205a29a6"

for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value() = $modelo[$r]
    if ($politica.ContainsKey($r)) {
        $ws.Cells.Item($r, 4).Value() = $politica[$r]
    }
    $ws.Cells.Item($r, 6).Value() = $tipo[$r]
    $ws.Cells.Item($r, 7).Value() = $link[$r]
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()